# "Generate Report for Handoff" — refresh the localization-status report:
#  - Priority column (E) on the zh-cn / de-de sheets flips from blank to
#    "ht" for the rows that just got a fresh handoff package generated.
#  - The "Latest Handoff Datetime" timestamps for those same rows (and the
#    matching "Latest HO Xliff Generate Date" on the Overview sheet) move
#    forward to the new generation time.

$wb = $excel.ActiveWorkbook

$rows = @(7, 9, 10, 11, 12, 14)

$zhcn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $zhcn.Range("E$r").Value = "ht"
    $zhcn.Range("H$r").Value = "2016-09-06 22:26:28"
}

$dede = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $dede.Range("E$r").Value = "ht"
    $dede.Range("H$r").Value = "2016-09-06 22:26:34"
}

$overview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $overview.Range("G$r").Value = "2016-09-06 22:26:34"
}
